$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.700.18"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.942.60"
$ws.Range("E3").Value = "  -2.50%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.50"
$ws.Range("E5").Value = "  -2.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.27"
$ws.Range("E6").Value = "  +1.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.936.77"
$ws.Range("E9").Value = "  -2.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.62"
$ws.Range("E10").Value = "  -2.46%  "
$ws.Range("E11").Value = "  -3.86%  "
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("E13").Value = "  -4.01%  "
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.672.59"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.432.25"
$ws.Range("E17").Value = "  -2.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.05"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.941.92"
$ws.Range("E19").Value = "  -2.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.83"
$ws.Range("E20").Value = "  +12.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "445.34"
$ws.Range("E21").Value = "  -2.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.695"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("E23").Value = "  -1.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.05"
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.23"
$ws.Range("E25").Value = "  -2.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.21"
$ws.Range("E26").Value = "  -1.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.99"
$ws.Range("E28").Value = "  -6.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.19"
$ws.Range("E29").Value = "  +1.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.45"
$ws.Range("E30").Value = "  +4.45%  "
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000100"
$ws.Range("E32").Value = "  -6.01%  "
$ws.Range("E33").Value = "  +4.90%  "
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.971"
$ws.Range("E36").Value = "  -2.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.69"
$ws.Range("E37").Value = "  -2.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "46.80"
$ws.Range("E38").Value = "  +7.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.08"
$ws.Range("E39").Value = "  -1.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.97"
$ws.Range("E40").Value = "  -9.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.301"
$ws.Range("E41").Value = "  -3.14%  "
$ws.Range("E42").Value = "  -1.47%  "
$ws.Range("E43").Value = "  -6.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.46"
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "380.50"
$ws.Range("E45").Value = "  -2.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0349"
$ws.Range("E46").Value = "  -1.76%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.665.92"
$ws.Range("E47").Value = "  -5.01%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.59"
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.86"
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.15"
$ws.Range("E51").Value = "  +0.57%  "
